$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New header cells, matching the header style used by the rest of row 1 (column H).
$ws.Range("I1").Value = "I0"
$ws.Range("J1").Value = "IF"
$ws.Range("H1").Copy()
$ws.Range("I1:J1").PasteSpecial(-4122)

# Data rows 2-41: I gets a constant 1, J mirrors column H's value for the row.
for ($r = 2; $r -le 41; $r++) {
    $hVal = $ws.Cells.Item($r, 8).Value2
    $ws.Cells.Item($r, 9).Value = 1
    $ws.Cells.Item($r, 10).Value = $hVal
}
